$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholder text (09/07/2013
#    -> 18/09/2014) on the slide master and every slide layout.
# ---------------------------------------------------------------------
$newDate = "18/09/2014"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes
for ($l = 1; $l -le $p.SlideMaster.CustomLayouts.Count; $l++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($l)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Group the logo artwork on slide 1 (the two ellipses, the two
#    highlight ellipses, the checkmark "Forma livre" shape and the
#    "FunTester" caption rectangle) into a single "Grupo 1" group.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$shapeNames = @("Elipse 3", "Elipse 4", "Elipse 5", "Elipse 6", "Forma livre 51", "Retângulo 52")
$idxs = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($shapeNames -contains $s.Shapes.Item($i).Name) {
        $idxs += $i
    }
}

$range = $s.Shapes.Range($idxs)
$group = $range.Group()
$group.Name = "Grupo 1"
